$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 5, mirroring row 4's layout with a new script line ---
$ws.Range("A5").Value = "SCRIPT/P01P04A/us3113.ssb"
$ws.Range("B5").Value = 83
$ws.Range("C5").Value = $ws.Range("C4").Value2
$ws.Range("D5").Value = $ws.Range("D4").Value2
$ws.Range("E5").Value = $ws.Range("E4").Value2

# D5:E5 keep the same formatting as D4:E4 (untouched style)
$ws.Range("D4:E4").Copy()
$ws.Range("D5:E5").PasteSpecial(-4122)   # xlPasteFormats

# Row 5 uses the same row height as row 4
$ws.Rows.Item(5).RowHeight = 43.2

# --- Box A4:C4 and A5:C5 with a thin top+bottom border ---
foreach ($r in 4, 5) {
    $rng = $ws.Range("A" + $r + ":C" + $r)

    $bottom = $rng.Borders.Item(9)
    $bottom.ColorIndex = 1
    $bottom.LineStyle = 1

    $top = $rng.Borders.Item(8)
    $top.ColorIndex = 1
    $top.LineStyle = 1
}

# --- Selection ends on E5 ---
[void]$ws.Range("E5").Select()

Write-Host "done"
